$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 490, pushing the existing rows 490-531 down to 492-533.
$ws.Range("A490:A491").EntireRow.Insert()

# Copy the formatting of the template rows (previously 490/491, now 492/493) into
# the freshly inserted rows so the date column keeps its date number format, etc.
$ws.Range("A492:R493").Copy()
$ws.Range("A490:R491").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 490: "Primera" quality entry for the new week.
$ws.Cells.Item(490, 1).Value = 8
$ws.Cells.Item(490, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(490, 3).Value = "Coquimbo"
$ws.Cells.Item(490, 4).Value = 44578
$ws.Cells.Item(490, 5).Value = 4
$ws.Cells.Item(490, 6).Value = 100112008
$ws.Cells.Item(490, 7).Value = "Coliflor"
$ws.Cells.Item(490, 8).Value = "Sin especificar"
$ws.Cells.Item(490, 9).Value = "Primera"
$ws.Cells.Item(490, 10).Value = 2700
$ws.Cells.Item(490, 11).Value = 650
$ws.Cells.Item(490, 12).Value = 700
$ws.Cells.Item(490, 13).Value = 675
$ws.Cells.Item(490, 14).Value = "`$/unidad"
$ws.Cells.Item(490, 15).Value = "Provincia del Elqu$([char]0x00ED)"
$ws.Cells.Item(490, 16).Value = 675
$ws.Cells.Item(490, 17).Value = 1
$ws.Cells.Item(490, 18).Value = "Hortaliza"

# Row 491: "Segunda" quality entry for the new week.
$ws.Cells.Item(491, 1).Value = 8
$ws.Cells.Item(491, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(491, 3).Value = "Coquimbo"
$ws.Cells.Item(491, 4).Value = 44578
$ws.Cells.Item(491, 5).Value = 4
$ws.Cells.Item(491, 6).Value = 100112008
$ws.Cells.Item(491, 7).Value = "Coliflor"
$ws.Cells.Item(491, 8).Value = "Sin especificar"
$ws.Cells.Item(491, 9).Value = "Segunda"
$ws.Cells.Item(491, 10).Value = 1600
$ws.Cells.Item(491, 11).Value = 550
$ws.Cells.Item(491, 12).Value = 600
$ws.Cells.Item(491, 13).Value = 575
$ws.Cells.Item(491, 14).Value = "`$/unidad"
$ws.Cells.Item(491, 15).Value = "Provincia del Elqu$([char]0x00ED)"
$ws.Cells.Item(491, 16).Value = 575
$ws.Cells.Item(491, 17).Value = 1
$ws.Cells.Item(491, 18).Value = "Hortaliza"
